$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Coin / Link / Price / Volume(1h) cells with the latest coinranking.com
# scrape. Three rank swaps are included: EthereumClassic<->Hedera (rows 32-33),
# Stellar<->RenderToken (rows 35-36), Celestia<->ARBITRUM (rows 44-45).
#
# Price cells are plain text in this sheet (even when the text looks like a
# number, e.g. "9.00" or "0.0360"), so numeric-looking values are written with a
# leading apostrophe - same as typing '9.00 into Excel - to keep Excel from
# silently re-typing them as numbers and dropping the significant trailing zeros.
$ws.Range("D2").Value = "43.836.06"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.314.71"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'114.02"
$ws.Range("E5").Value = "  +19.94%  "
$ws.Range("D6").Value = "'272.08"
$ws.Range("E6").Value = "  +1.86%  "
$ws.Range("E7").Value = "  +0.54%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("D9").Value = "'0.629"
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("D10").Value = "'47.58"
$ws.Range("E10").Value = "  +8.03%  "
$ws.Range("D11").Value = "'0.0946"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "'9.00"
$ws.Range("E12").Value = "  +15.39%  "
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").Value = "'15.87"
$ws.Range("E14").Value = "  +4.16%  "
$ws.Range("D15").Value = "2.655.75"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "'0.866"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "2.303.64"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "43.838.85"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "'6.72"
$ws.Range("E20").Value = "  +7.86%  "
$ws.Range("D21").Value = "'72.69"
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'2.51"
$ws.Range("E22").Value = "  +6.76%  "
$ws.Range("D23").Value = "'234.89"
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("E24").Value = "  +7.14%  "
$ws.Range("D25").Value = "'2.90"
$ws.Range("E25").Value = "  +16.04%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'11.59"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "'43.52"
$ws.Range("E28").Value = "  +16.24%  "
$ws.Range("D29").Value = "'3.44"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'177.97"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'21.89"
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0939"
$ws.Range("E33").Value = "  +6.02%  "
$ws.Range("D34").Value = "'5.65"
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.127"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'4.81"
$ws.Range("E36").Value = "  +8.58%  "
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("D38").Value = "'3.97"
$ws.Range("E38").Value = "  +21.41%  "
$ws.Range("D39").Value = "'0.0360"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("D40").Value = "'0.246"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "'2.41"
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").Value = "'69.39"
$ws.Range("E42").Value = "  +11.99%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").Value = "'12.68"
$ws.Range("E44").Value = "  +6.76%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").Value = "'1.38"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").Value = "'5.72"
$ws.Range("E46").Value = "  +9.21%  "
$ws.Range("D47").Value = "'8.86"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'100.39"
$ws.Range("E49").Value = "  +1.27%  "
$ws.Range("D50").Value = "'1.23"
$ws.Range("E50").Value = "  +3.29%  "
$ws.Range("D51").Value = "'0.463"
$ws.Range("E51").Value = "  +10.19%  "
